$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 750
$ws1.Range("F4").Value = 1493
$ws1.Range("F5").Value = 228
$ws1.Range("F6").Value = 94
$ws1.Range("F8").Value = 6228
$ws1.Range("F9").Value = 70
$ws1.Range("F10").Value = 404
$ws1.Range("F12").Value = 5176
$ws1.Range("F13").Value = 27
$ws1.Range("F14").Value = 179
$ws1.Range("F15").Value = 1178
$ws1.Range("F17").Value = 362
$ws1.Range("F18").Value = 65
$ws1.Range("F19").Value = 10
$ws1.Range("F20").Value = 297
$ws1.Range("F22").Value = 3696
$ws1.Range("F23").Value = 152

# Sheet "演出" (Performances) - sheet2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 80

# Sheet "全部类型" (All types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 80
$ws4.Range("F4").Value = 750
$ws4.Range("F5").Value = 1493
$ws4.Range("F6").Value = 228
$ws4.Range("F7").Value = 94
$ws4.Range("F9").Value = 6228
$ws4.Range("F10").Value = 70
$ws4.Range("F11").Value = 404
$ws4.Range("F13").Value = 5176
$ws4.Range("F14").Value = 27
$ws4.Range("F15").Value = 179
$ws4.Range("F16").Value = 1178
$ws4.Range("F18").Value = 362
$ws4.Range("F19").Value = 65
$ws4.Range("F20").Value = 10
$ws4.Range("F21").Value = 297
$ws4.Range("F23").Value = 3696
$ws4.Range("F25").Value = 152
